# Correction in SA algorithm and 746 logs
# Update Fitness (column C) values:
#   Rows 2-116  -> 7310
#   Rows 117-252 -> 7293

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C116").Value = 7310
$ws.Range("C117:C252").Value = 7293
